$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.114.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.06%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.813.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.81%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.73%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4623"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.20%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3756"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07416"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8620"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.70%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.61"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.42%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.816.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.654"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.387"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07090"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.92"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.22%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.14%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008738"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.04%  "

$ws.Range("E19").Value = "  +0.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.96%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.110.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.98%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.316"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.37%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.61%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.047.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.99%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.923"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.51%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.213"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.99%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.279"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.81%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.22%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08927"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.85%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7727"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.69%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.170"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.533"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.97%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.886"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("E36").Value = "  +0.11%  "

$ws.Range("E37").Value = "  +3.50%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01960"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05229"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.95%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.250"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.06%  "

$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.921"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.374"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +19.23%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5282"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1678"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.23%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.610"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5041"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.86%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.000"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.671"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.63%  "

$ws.Range("E51").Value = "  +0.10%  "
